# Applies the term-list / body-part tagging edit described by the diff:
#   - wraps "front" (split across runs as "fro" + <exp>n</exp> + "t") in <bp>…</bp>
#     inside the <head> of div p033v_1
#   - wraps "gecton" and "piece d'<m>argent</m>" in <tl>…</tl> tags inside the
#     first <ab> of div p033v_1
#
# Each target substring is the *entire* text of a single run, so the edits are
# scoped to the exact paragraph (found by locating the unique source text) and
# applied either via Find/Replace (whole-run, unique text) or by locating a
# tiny Range via Find and assigning .Text directly (for the 1-character "t"
# run, where a plain Find for "t" would be ambiguous).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "<head>… au fro<exp>n</exp>t</head>"  ->  "… au <bp>fro<exp>n</exp>t</bp></head>"
# ---------------------------------------------------------------------

# Locate the paragraph that holds " au fro" (the <head> continuation line).
$headPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*piece d*au fro*") {
        $headPara = $cand
        break
    }
}

# 1a) " au fro"  ->  " au <bp>fro"   (single run, unique in the paragraph)
$r1 = $headPara.Range
$r1.Find.Execute(" au fro", $true, $false, $false, $false, $false, $true, 1, $false, " au <bp>fro", 2)

# 1b) the lone "t" run right after "<exp>n</exp>" and before "</head>"
#     -> "t</bp>"
$r2 = $headPara.Range
$r2.Find.Execute("</exp>")
$tRange = $d.Range($r2.End, $r2.End + 1)
$tRange.Text = "t</bp>"

# ---------------------------------------------------------------------
# 2) "<ab>Ayes un gecton ou piece d'<m>argent</m> et la mouille…"
#    -> "<ab>Ayes un <tl>gecton</tl> ou <tl>piece d'<m>argent</m></tl> et la mouille…"
# ---------------------------------------------------------------------

# Locate the paragraph that holds "Ayes un gecton".
$abPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Ayes un gecton*") {
        $abPara = $cand
        break
    }
}

# 2a) "Ayes un gecton ou piece d" -> "Ayes un <tl>gecton</tl> ou <tl>piece d"
$r3 = $abPara.Range
$r3.Find.Execute("Ayes un gecton ou piece d", $true, $false, $false, $false, $false, $true, 1, $false, "Ayes un <tl>gecton</tl> ou <tl>piece d", 2)

# 2b) "</m>" (closing the <m>argent</m> right after "piece d'") -> "</m></tl>"
$r4 = $abPara.Range
$r4.Find.Execute("</m>", $true, $false, $false, $false, $false, $true, 1, $false, "</m></tl>", 2)
